# Fruta / hortaliza, semanal
#
# A new weekly price observation for "Apio" (Macroferia Regional de Talca)
# is inserted as a new data row right after the header block of existing
# rows, at worksheet row 89 (the sheet's data rows run from row 2 to the
# end). Inserting the row shifts every existing row from 89 downward by
# one (old row 89 -> new row 90, ..., old row 112 -> new row 113), which
# also grows the sheet's used range from R112 to R113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 89..end down by one, opening up a blank row 89.
$ws.Rows.Item(89).Insert()

# Populate the newly-opened row 89 with the new weekly record.
$ws.Range("A89").Value = 5
$ws.Range("B89").Value = "Macroferia Regional de Talca"
$ws.Range("C89").Value = "Maule"
$ws.Range("D89").Value = 44463
$ws.Range("E89").Value = 7
$ws.Range("F89").Value = 100112017
$ws.Range("G89").Value = "Apio"
$ws.Range("H89").Value = "Americana (o)"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 8000
$ws.Range("M89").Value = 8000
$ws.Range("N89").Value = '$/docena de matas'
$ws.Range("O89").Value = "Provincia del Elquí"
$ws.Range("P89").Value = 1333
$ws.Range("Q89").Value = 6
$ws.Range("R89").Value = "Hortaliza"
